$p = $ppt.ActivePresentation

# --- Slide 5 ("5_Publish"): content placeholder has four "Help ..." paragraphs
# that need "Help" -> "Helps" with the run split into three runs:
#   "Helps" | " " | "<rest of original first run text>"
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$tr5 = $sh5.TextFrame.TextRange

$paraIndexes = @(1, 3, 5, 6)
foreach ($idx in $paraIndexes) {
    $para = $tr5.Paragraphs($idx)
    $c1 = $para.Characters(1, 4)
    $c1.Text = "Helps"
    $para2 = $tr5.Paragraphs($idx)
    $c2 = $para2.Characters(6, 1)
    $c2.Text = " "
}

# --- Slide 7: single paragraph where "Help " is already its own run,
# just needs the text updated to "Helps ".
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(2)
$tr7 = $sh7.TextFrame.TextRange
$para7 = $tr7.Paragraphs(4)
$c7 = $para7.Characters(1, 5)
$c7.Text = "Helps "
